# Apply the "Feedback" + "User Login" test-case content to the document,
# move the _GoBack bookmark to a new trailing empty paragraph, and change
# the page size/margins (US Letter -> A4-ish 11900x16840 twips, left/right
# margins 1440 -> 1800 twips).

$d = $word.ActiveDocument

# ------------------------------------------------------------------
# 1. The last paragraph of the "Notifications" section currently ends
#    with the `_GoBack` bookmark inside it. Remove it from there; it
#    will be re-added to a new trailing empty paragraph later.
# ------------------------------------------------------------------
if ($d.Bookmarks.Exists("_GoBack")) {
    $d.Bookmarks.Item("_GoBack").Delete()
}

# ------------------------------------------------------------------
# 2. Remove the old "Feedback" section (heading + single test case +
#    trailing blank paragraphs) that followed the Notifications section.
#    It is being replaced wholesale by the much larger block below.
# ------------------------------------------------------------------
$anchor = $d.Paragraphs.Item(5)
if ($d.Paragraphs.Count -gt 5) {
    $startP = $d.Paragraphs.Item(6)
    $endP = $d.Paragraphs.Item($d.Paragraphs.Count)
    $oldRange = $d.Range($startP.Range.Start, $endP.Range.End)
    $oldRange.Delete()
}

# ------------------------------------------------------------------
# 3. Insert the new "Feedback" and "User Login" sections paragraph by
#    paragraph, right after the anchor paragraph.
# ------------------------------------------------------------------
$newParas = @(
  @{ Style = "Heading1"; Text = "Feedback" },
  @{ Style = "Heading2"; Text = "Test Case 1: Login with correct data" },
  @{ Style = "Normal"; Text = "Precondition: The correct user details are entered, the user login button is pressed." },
  @{ Style = "Normal"; Text = "Post-condition: The login was a success. The user is navigated to the home page." },
  @{ Style = "Normal"; Text = "Expected result: Successful  login accompanied by a message from the website indicating that the user  logged in successfully." },
  @{ Style = "Normal"; Text = "" },
  @{ Style = "Heading2"; Text = "Test Case 2: Login with incorrect data" },
  @{ Style = "Normal"; Text = "Precondition: Incorrect user details are used in order to login." },
  @{ Style = "Normal"; Text = "Post-condition: Login was a not a success. A message is displayed to the user that their login details are not correct. " },
  @{ Style = "Normal"; Text = "Expected result: Unsuccessful login accompanied with a message to indicate why." },
  @{ Style = "Normal"; Text = "" },
  @{ Style = "Heading2"; Text = "Test Case 3: Register with correct data" },
  @{ Style = "Normal"; Text = "Precondition: Correct data is given in order to successfully register." },
  @{ Style = "Normal"; Text = "Post-condition: Registration was a success. A message is displayed to the user to indicate that they have successfully registered and will now be able to log into the application. " },
  @{ Style = "Normal"; Text = "Expected result: Successful registration accompanied with a message to indicate this." },
  @{ Style = "Normal"; Text = "" },
  @{ Style = "Heading2"; Text = "Test Case 4: Register with incorrect data" },
  @{ Style = "Normal"; Text = "Precondition: Incorrect data is used to register." },
  @{ Style = "Normal"; Text = "Post-condition: Register button greys out when incorrect data is given and user is not able to register until correct data is entered." },
  @{ Style = "Normal"; Text = "Expected result: Registration not successful accompanied by a message or feedback to indicate why." },
  @{ Style = "Heading2"; Text = "Test Case 5: Update Account details" },
  @{ Style = "Normal"; Text = "Precondition: New data is entered to update account details." },
  @{ Style = "Normal"; Text = "Post-condition: The user details are updated however no message is displayed to inform the user." },
  @{ Style = "Normal"; Text = "Expected result:  A notification should be displayed to indicate that the account details were successfully updated." },
  @{ Style = "Heading2"; Text = "Test Case 6: Save route" },
  @{ Style = "Normal"; Text = "Precondition: A user requests to save a route." },
  @{ Style = "Normal"; Text = "Post-condition: The route is saved. No feedback is given." },
  @{ Style = "Normal"; Text = "Expected result: A message displayed to the user to indicate that the route was successfully saved." },
  @{ Style = "Heading2"; Text = "Test Case 7: Save point of interest" },
  @{ Style = "Normal"; Text = "Precondition: A user requests to save a point of interest." },
  @{ Style = "Normal"; Text = "Post-condition: The point of interest is saved and a message is displayed to the user to indicate this." },
  @{ Style = "Normal"; Text = "Expected result: A message should be displayed to indicate that the point of interest was successfully saved." },
  @{ Style = "Heading2"; Text = "Test Case 7: Log out" },
  @{ Style = "Normal"; Text = "Precondition: A user requests to log out of the application." },
  @{ Style = "Normal"; Text = "Post-condition: The user is logged out of the website and returned to the login page, no message is displayed." },
  @{ Style = "Normal"; Text = "Expected result:  A message should be displayed to indicate that the user successfully logged out of the website." },
  @{ Style = "Normal"; Text = "" },
  @{ Style = "Heading2"; Text = "Test Case 8: Log out" },
  @{ Style = "Normal"; Text = "Precondition: A user requests to log out of the application." },
  @{ Style = "Normal"; Text = "Post-condition: The user is logged out of the website and returned to the login page, no message is displayed." },
  @{ Style = "Normal"; Text = "Expected result: A message to display that the user successfully logged out of the website." },
  @{ Style = "Heading2"; Text = "Test Case 9: Delete GIS object" },
  @{ Style = "Normal"; Text = "Precondition: An admin requests to delete a GIS object." },
  @{ Style = "Normal"; Text = "Post-condition: The GIS object is deleted and a message is displayed to confirm it was deleted successfully." },
  @{ Style = "Normal"; Text = "Expected result: A message is displayed to indicate that the GIS object was deleted." },
  @{ Style = "Heading2"; Text = "Test Case 9: Save changes to GIS object" },
  @{ Style = "Normal"; Text = "Precondition: An admin makes changes to a GIS object and requests to save these changes." },
  @{ Style = "Normal"; Text = "Post-condition: The changes to the GIS object is made, no message is displayed." },
  @{ Style = "Normal"; Text = "Expected result: A notification is displayed to indicate the changes were saved." },
  @{ Style = "Normal"; Text = "" },
  @{ Style = "Heading2"; Text = "Test Case 11: Add location" },
  @{ Style = "Normal"; Text = "Precondition: An admin wants to add a location." },
  @{ Style = "Normal"; Text = "Post-condition: The location is added successfully. A notification is displayed to confirm this." },
  @{ Style = "Normal"; Text = "Expected result: A notification is displayed to inform the user that the location was successfully added." },
  @{ Style = "Normal"; Text = "" },
  @{ Style = "Heading2"; Text = "Test Case 12: Remove a location" },
  @{ Style = "Normal"; Text = "Precondition: An admin requests to remove a location." },
  @{ Style = "Normal"; Text = "Post-condition: The location is removed and a message is displayed to confirm this." },
  @{ Style = "Normal"; Text = "Expected result: A notification is displayed to show that the location was successfully removed." },
  @{ Style = "Normal"; Text = "" },
  @{ Style = "Heading2"; Text = "Test Case 13: Remove a location" },
  @{ Style = "Normal"; Text = "Precondition: An admin requests to remove a location." },
  @{ Style = "Normal"; Text = "Post-condition: The location is removed and a message is displayed to confirm this." },
  @{ Style = "Normal"; Text = "Expected result: A notification is displayed to show that the location was successfully removed." },
  @{ Style = "Normal"; Text = "" },
  @{ Style = "Heading2"; Text = "Test Case 14: Add and remove admin rights" },
  @{ Style = "Normal"; Text = "Precondition: An admin requests to add or remove the admin rights of a user." },
  @{ Style = "Normal"; Text = "Post-condition: Admin rights are added or removed from a user, no notification or other feedback is given to confirm this." },
  @{ Style = "Normal"; Text = "Expected result: A notification is displayed when admin rights are added or removed from a user." },
  @{ Style = "Normal"; Text = "" },
  @{ Style = "Heading2"; Text = "Test Case 15: Remove user" },
  @{ Style = "Normal"; Text = "Precondition: An admin requests to remove a user." },
  @{ Style = "Normal"; Text = "Post-condition: The user is removed, no notification is displayed." },
  @{ Style = "Normal"; Text = "Expected result: A notification to be displayed when a user is removed." },
  @{ Style = "Heading1"; Text = "User Login" },
  @{ Style = "Heading2"; Text = "Test Case 1: Login with correct data" },
  @{ Style = "Normal"; Text = "Precondition: The correct user details are entered." },
  @{ Style = "Normal"; Text = "Post-condition: Login was a success." },
  @{ Style = "Normal"; Text = "Expected result: Successful login." },
  @{ Style = "Normal"; Text = "" },
  @{ Style = "Heading2"; Text = "Test Case 2: Login with no data entered " },
  @{ Style = "Normal"; Text = "Precondition: No data is entered into the required login input boxes; the user login button is pressed." },
  @{ Style = "Normal"; Text = "Post-condition: Login was not successful." },
  @{ Style = "Normal"; Text = "Expected result: Login not successful." },
  @{ Style = "Normal"; Text = "" },
  @{ Style = "Heading2"; Text = "Test Case 3: Login with incorrect password " },
  @{ Style = "Normal"; Text = "Precondition: A correct username is entered with an incorrect password." },
  @{ Style = "Normal"; Text = "Post-condition: Login was not successful." },
  @{ Style = "Normal"; Text = "Expected result: Login not successful." },
  @{ Style = "Normal"; Text = "" },
  @{ Style = "Heading2"; Text = "Test Case 4: Login with incorrect username and password" },
  @{ Style = "Normal"; Text = "Precondition: An incorrect username and password are entered." },
  @{ Style = "Normal"; Text = "Post-condition: Login was not successful." },
  @{ Style = "Normal"; Text = "Expected result: Login not successful." },
  @{ Style = "Normal"; Text = "" },
  @{ Style = "Normal"; Text = "" },
  @{ Style = "Normal"; Text = "" },
  @{ Style = "Normal"; Text = "" },
  @{ Style = "Normal"; Text = "" }
)

$cur = $anchor
foreach ($item in $newParas) {
    $cur.Range.InsertParagraphAfter()
    $newp = $d.Paragraphs.Item($cur.Index + 1)
    $newp.Style = $item.Style
    if ($item.Text -ne "") {
        $newp.Range.Text = $item.Text
    }
    $cur = $newp
}

# ------------------------------------------------------------------
# 4. Re-add the `_GoBack` bookmark as its own empty paragraph at the
#    very end of the document (after the last blank paragraph).
# ------------------------------------------------------------------
$cur.Range.InsertParagraphAfter()
$bookmarkPara = $d.Paragraphs.Item($cur.Index + 1)
$bookmarkPara.Style = "Normal"
$d.Bookmarks.Add("_GoBack", $bookmarkPara.Range)

# ------------------------------------------------------------------
# 5. Update the section's page size / margins.
# ------------------------------------------------------------------
$ps = $d.PageSetup
$ps.PageWidth = 595
$ps.PageHeight = 842
$ps.LeftMargin = 90
$ps.RightMargin = 90

Write-Output ("Paragraphs: " + $d.Paragraphs.Count)
